$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "lunges pad"
$ws.Cells.Item(2, 1).Value = "knee protection pad"
$ws.Cells.Item(3, 1).Value = "boys basketball pants"
$ws.Cells.Item(4, 1).Value = "working knee pads for men"
$ws.Cells.Item(5, 1).Value = "knee yoga pants"
$ws.Cells.Item(6, 1).Value = "padded knee sleeve for sliding"
$ws.Cells.Item(7, 1).Value = "men tight pants"
$ws.Cells.Item(8, 1).Value = "protective compression wear"
$ws.Cells.Item(9, 1).Value = "spandex compression shorts men"
$ws.Cells.Item(10, 1).Value = "softball mens pants"
$ws.Cells.Item(11, 1).Value = "compression capri"
$ws.Cells.Item(12, 1).Value = "softball sliding pants youth girls"
$ws.Cells.Item(13, 1).Value = "baseball pants black mens"
$ws.Cells.Item(14, 1).Value = "wrestling knee pads pair"
$ws.Cells.Item(15, 1).Value = "compressions tights for men"
$ws.Cells.Item(16, 1).Value = "men sports leggings"
$ws.Cells.Item(17, 1).Value = "compression shorts boys padded"
$ws.Cells.Item(18, 1).Value = "basketball padded knee sleeve"
$ws.Cells.Item(19, 1).Value = "yoga pad thick"
$ws.Cells.Item(20, 1).Value = "mens leggings tall"
$ws.Cells.Item(21, 1).Value = "knee protector construction"
$ws.Cells.Item(22, 1).Value = "basketball leggings for girls"
$ws.Cells.Item(23, 1).Value = "lacrosse compression shorts"
$ws.Cells.Item(24, 1).Value = "boys sports tights leggings"
$ws.Cells.Item(25, 1).Value = "mens sweat pads"
$ws.Cells.Item(26, 1).Value = "mens construction knee pads"
$ws.Cells.Item(27, 1).Value = "yoga pants men"
$ws.Cells.Item(28, 1).Value = "bursitis knee"
$ws.Cells.Item(29, 1).Value = "thread protector paintball"
$ws.Cells.Item(30, 1).Value = "tight gym pants men"
$ws.Cells.Item(31, 1).Value = "adult bee tights"
$ws.Cells.Item(32, 1).Value = "basketball shorts in bulk"
$ws.Cells.Item(33, 1).Value = "rash guard men leggings"
$ws.Cells.Item(34, 1).Value = "youth sliding short"
$ws.Cells.Item(35, 1).Value = "running knee compression"
$ws.Cells.Item(36, 1).Value = "knee protection for running"
$ws.Cells.Item(37, 1).Value = "girls basketball leggings"
$ws.Cells.Item(38, 1).Value = "youth xxl baseball pants"
$ws.Cells.Item(39, 1).Value = "boys softball pants"
$ws.Cells.Item(40, 1).Value = "tight for boys"
$ws.Cells.Item(41, 1).Value = "big and tall leggings men"
$ws.Cells.Item(42, 1).Value = "knee protect"
$ws.Cells.Item(43, 1).Value = "knee sleeve basketball men"
$ws.Cells.Item(44, 1).Value = "protective baseball"
$ws.Cells.Item(45, 1).Value = "compression pants for recovery"
$ws.Cells.Item(46, 1).Value = "girls basketball knee guards"
$ws.Cells.Item(47, 1).Value = "basketball knee pads for women"
$ws.Cells.Item(48, 1).Value = "gel wrestling knee pads"
$ws.Cells.Item(49, 1).Value = "mens hockey pads"
$ws.Cells.Item(50, 1).Value = "long knee pads volleyball"
$ws.Cells.Item(51, 1).Value = "mens compression pants 3/4 length"
$ws.Cells.Item(52, 1).Value = "knee pads motorcycle"
$ws.Cells.Item(53, 1).Value = "sliding shorts youth girls"
$ws.Cells.Item(54, 1).Value = "boys athletic leggings youth"
$ws.Cells.Item(55, 1).Value = "knee pads for basketball women"
$ws.Cells.Item(56, 1).Value = "big mens compression pants"
$ws.Cells.Item(57, 1).Value = "soccer sliding pants"
$ws.Cells.Item(58, 1).Value = "good thread mens pants"
$ws.Cells.Item(59, 1).Value = "xl knee pads for men"
$ws.Cells.Item(60, 1).Value = "capri pouches adults"
$ws.Cells.Item(61, 1).Value = "knee sleeves hex"
$ws.Cells.Item(62, 1).Value = "knee compression sleeve with padding"
$ws.Cells.Item(63, 1).Value = "baleaf mens pants"
$ws.Cells.Item(64, 1).Value = "leggings tight"
$ws.Cells.Item(65, 1).Value = "compression tights mens"
$ws.Cells.Item(66, 1).Value = "compression knee sleeve basketball"
$ws.Cells.Item(67, 1).Value = "volleyball knee pads gel"
$ws.Cells.Item(68, 1).Value = "knee pads for work xxl"
$ws.Cells.Item(69, 1).Value = "knee pad for exercise"
$ws.Cells.Item(70, 1).Value = "running tights youth boys"
$ws.Cells.Item(71, 1).Value = "boy capri pants"
$ws.Cells.Item(72, 1).Value = "youth hockey girdle"
$ws.Cells.Item(73, 1).Value = "youth compression leggings boys"
$ws.Cells.Item(74, 1).Value = "sliding compression shorts"
$ws.Cells.Item(75, 1).Value = "baseball leg protection"
$ws.Cells.Item(76, 1).Value = "knee sleeves for basketball youth"
$ws.Cells.Item(77, 1).Value = "youth paintball pants"
$ws.Cells.Item(78, 1).Value = "youth boy tights"
$ws.Cells.Item(79, 1).Value = "softball catcher pants"
$ws.Cells.Item(80, 1).Value = "knee pads volleyball mens"
$ws.Cells.Item(81, 1).Value = "mens running compression pants"
$ws.Cells.Item(82, 1).Value = "football tights youth boys"
$ws.Cells.Item(83, 1).Value = "athletic tights youth boys"
$ws.Cells.Item(84, 1).Value = "padded calf sleeve"
$ws.Cells.Item(85, 1).Value = "mens softball pants"
$ws.Cells.Item(86, 1).Value = "mens baseball shorts"
$ws.Cells.Item(87, 1).Value = "knee pad for soccer"
$ws.Cells.Item(88, 1).Value = "compression leggings boys"
$ws.Cells.Item(89, 1).Value = "capri for men"
$ws.Cells.Item(90, 1).Value = "mens workout pants leggings"
$ws.Cells.Item(91, 1).Value = "knee pad sport"
$ws.Cells.Item(92, 1).Value = "tights pants boys"
$ws.Cells.Item(93, 1).Value = "softball sliding shorts women"
$ws.Cells.Item(94, 1).Value = "padded sliding shorts"
$ws.Cells.Item(95, 1).Value = "arthritis hope knee sleeve"
$ws.Cells.Item(96, 1).Value = "padded knee sleeve basketball"
$ws.Cells.Item(97, 1).Value = "snowboarding mens pants"
$ws.Cells.Item(98, 1).Value = "knee compression cold"
$ws.Cells.Item(99, 1).Value = "knee pads squats"
$ws.Cells.Item(100, 1).Value = "men work pants with knee pads"
